# Hindalco price sheet daily update: shift all existing rows down by one
# and insert a brand new "today" row (31-10-2025) at the top.
#
# The source table keeps the most-recent date in row 2 and pushes every
# earlier row down by one each day a new price is published; the oldest
# row that falls off the bottom of the previously-used range re-appears
# as a brand-new last row (since its data simply moves down with the
# rest of the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top of the data (row 2, right below the header).
# This shifts all data rows down by one (cell values only -- hyperlinks
# are rebuilt explicitly below), which automatically produces the new
# row 143 (a copy of what used to be the last row, 142).
$ws.Rows("2:2").Insert()

# The newly inserted row 2 is blank; populate it by duplicating row 3
# (which now holds what used to be row 2's data) and then overwrite the
# date with the new day's date. Values and formats are pasted separately
# so the data rows' cell styles (s=3 text / s=4 price) are reproduced
# exactly instead of Excel synthesizing a new blended style.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A2").Value = "31-10-2025"

# Row-insertion only moves cell contents; it leaves the worksheet's
# Hyperlinks collection anchored to their original (pre-shift) cells, so
# every hyperlink is now one row off from the circular-link text it
# should decorate (and the row that used to be empty, now row 82, needs
# a brand new one). Rebuild the whole collection from scratch: drop all
# existing hyperlinks, then re-add one per non-blank Circular Link cell
# pointing at that cell's own text.
$ws.Hyperlinks.Delete()
for ($r = 2; $r -le 143; $r++) {
    $cell = $ws.Range("F$r")
    $url = $cell.Value()
    if ($url -ne "" -and $url -ne $null) {
        $ws.Hyperlinks.Add($cell, $url)
        # Adding a hyperlink auto-applies Excel's built-in blue/underline
        # "Hyperlink" style, but the source data keeps these cells in the
        # same plain style as the rest of the row (s=3). Re-stamp the
        # plain format (borrowed from that row's own A cell, which is
        # never touched by the hyperlink) over the link cell afterwards.
        $ws.Range("A$r").Copy()
        $cell.PasteSpecial(-4122)  # xlPasteFormats
    }
}
$excel.CutCopyMode = 0
